{"js": "// Update the session Date / Start Time / End Time / Total Time lines at the\n// top of the checklist.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\nconst startParagraph = paragraphs.items[1];\nconst endParagraph = paragraphs.items[2];\nconst totalParagraph = paragraphs.items[3];\n\n// Date: 21 September 2023 -> Date: 23 September 2023\nconst dateDigit = dateParagraph.search(\"1\", { matchCase: true });\ndateDigit.load(\"items\");\nawait context.sync();\ndateDigit.items[0].insertText(\"3\", Word.InsertLocation.replace);\n\n// Start Time: 6:00 PM -> Start Time: 12:46 PM\nconst startValue = startParagraph.search(\"6:00\", { matchCase: true });\nstartValue.load(\"items\");\nawait context.sync();\nstartValue.items[0].insertText(\"12:46\", Word.InsertLocation.replace);\n\n// End Time: 8:00 PM -> End Time:  6:46 PM\nconst endHour = endParagraph.search(\"8\", { matchCase: true });\nendHour.load(\"items\");\nawait context.sync();\nendHour.items[0].insertText(\" 6:\", Word.InsertLocation.replace);\n\nconst endMinutes = endParagraph.search(\":00 PM\", { matchCase: true });\nendMinutes.load(\"items\");\nawait context.sync();\nendMinutes.items[0].insertText(\"46 PM\", Word.InsertLocation.replace);\n\n// Total Time: 12 hour 23 Minutes -> Total Time: 9 hour 23 Minutes\nconst totalValue = totalParagraph.search(\"12\", { matchCase: true });\ntotalValue.load(\"items\");\nawait context.sync();\ntotalValue.items[0].insertText(\"9\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Update the session Date / Start Time / End Time / Total Time lines at the\n# top of the checklist.\n$d = $word.ActiveDocument\n\n# Date: 21 September 2023 -> Date: 23 September 2023\n$dateParagraph = $d.Paragraphs(1).Range\n$dateFind = $dateParagraph.Find\n$dateFind.Text = \"1\"\n$dateFind.Execute() | Out-Null\n$dateParagraph.Text = \"3\"\n\n# Start Time: 6:00 PM -> Start Time: 12:46 PM\n$startParagraph = $d.Paragraphs(2).Range\n$startFind = $startParagraph.Find\n$startFind.Text = \"6:00\"\n$startFind.Execute() | Out-Null\n$startParagraph.Text = \"12:46\"\n\n# End Time: 8:00 PM -> End Time:  6:46 PM\n$endParagraph = $d.Paragraphs(3).Range\n$endHourFind = $endParagraph.Find\n$endHourFind.Text = \"8\"\n$endHourFind.Execute() | Out-Null\n$endParagraph.Text = \" 6:\"\n\n$endMinutesRange = $d.Paragraphs(3).Range\n$endMinutesFind = $endMinutesRange.Find\n$endMinutesFind.Text = \":00 PM\"\n$endMinutesFind.Execute() | Out-Null\n$endMinutesRange.Text = \"46 PM\"\n\n# Total Time: 12 hour 23 Minutes -> Total Time: 9 hour 23 Minutes\n$totalParagraph = $d.Paragraphs(4).Range\n$totalFind = $totalParagraph.Find\n$totalFind.Text = \"12\"\n$totalFind.Execute() | Out-Null\n$totalParagraph.Text = \"9\"\n"}
